$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.804605722427368
$ws.Range("B1").Value = 6.182538032531738
$ws.Range("C1").Value = 5.307006359100342
$ws.Range("D1").Value = 6.167904376983643
$ws.Range("E1").Value = 3.7728111743927
